$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$results = @{
    2  = "Good"
    3  = "No"
    4  = "Good"
    5  = "No"
    6  = "No"
    7  = "Good"
    8  = "Good"
    9  = "No"
    10 = "No"
    11 = "Good"
    12 = "No"
    13 = "Good"
    14 = "Good"
    15 = "No"
    16 = "No"
    17 = "No"
    18 = "No"
    19 = "Good"
    20 = "Good"
    21 = "Good"
    22 = "Good"
    23 = "Good"
    24 = "No"
    25 = "Good"
    26 = "No"
    27 = "Good"
}

foreach ($row in $results.Keys) {
    $ws.Range("F$row").Value = $results[$row]
}
